$d = $word.ActiveDocument

$replacements = @(
    @("2024-10-31 Thursday", "2024-11-01 Friday"),
    @("35×28=980", "18×30=540"),
    @("91×41=3731", "59×74=4366"),
    @("48×78=3744", "22×42=924"),
    @("17×34=578", "87×30=2610"),
    @("43×15=645", "83×95=7885"),
    @("96×82=7872", "57×89=5073"),
    @("84×84=7056", "97×94=9118"),
    @("97×82=7954", "25×51=1275"),
    @("88×76=6688", "79×55=4345"),
    @("84×28=2352", "49×77=3773"),
    @("37×39=1443", "66×74=4884"),
    @("88×95=8360", "31×69=2139"),
    @("22×78=1716", "30×47=1410"),
    @("57×62=3534", "39×95=3705"),
    @("30×92=2760", "99×58=5742"),
    @("16×22=352", "45×59=2655"),
    @("67×74=4958", "81×98=7938"),
    @("82×33=2706", "39×36=1404"),
    @("55×38=2090", "65×49=3185"),
    @("67×38=2546", "52×50=2600"),
    @("13×97=1261", "67×80=5360"),
    @("22×15=330", "55×84=4620"),
    @("18×89=1602", "11×34=374"),
    @("97×48=4656", "38×99=3762"),
    @("20×69=1380", "65×32=2080")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
